# repull data, push all data, mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" column (F) values following a data repull
$ws.Range("F2").Value = -3
$ws.Range("F3").Value = -1
$ws.Range("F7").Value = -9
$ws.Range("F11").Value = -12
$ws.Range("F16").Value = -2
